# RRPONTSYD.xlsx update
# - Append 10 new daily observations (2023-08-09 .. 2023-08-22) to the "Data" sheet.
# - Refresh the FRED series metadata on "SeriesInfo" to match the new pull.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

# --- Data sheet: new rows 463-472 ---------------------------------------
# Copy the date-column formatting (style applied to A462) down across the
# new rows first, then fill in the real date-serial / value pairs.
$dataSheet.Range("A462").Copy($dataSheet.Range("A463:A472"))

$newRows = @(
    @{ Row = 463; Date = 45147; Value = 1796.519 },
    @{ Row = 464; Date = 45148; Value = 1759.897 },
    @{ Row = 465; Date = 45149; Value = 1773.236 },
    @{ Row = 466; Date = 45152; Value = 1799.311 },
    @{ Row = 467; Date = 45153; Value = 1743.784 },
    @{ Row = 468; Date = 45154; Value = 1796.725 },
    @{ Row = 469; Date = 45155; Value = 1794.120 },
    @{ Row = 470; Date = 45156; Value = 1819.201 },
    @{ Row = 471; Date = 45159; Value = 1824.788 },
    @{ Row = 472; Date = 45160; Value = 1812.294 }
)

foreach ($r in $newRows) {
    $dataSheet.Cells.Item($r.Row, 1).Value = $r.Date
    $dataSheet.Cells.Item($r.Row, 2).Value = $r.Value
}

# --- SeriesInfo sheet: refreshed metadata from FRED ----------------------
# B3/B4/B7 hold plain "YYYY-MM-DD" text (not real dates) in the source data.
# Force text formatting for the write so Excel doesn't auto-convert the
# literal string into a date serial, then drop the formatting override again
# so the cell's style stays exactly as it was before.
foreach ($addr in @("B3", "B4", "B7")) {
    $cell = $infoSheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "2023-08-22"
    $cell.ClearFormats()
}

$infoSheet.Range("B14").Value = "2023-08-22 13:01:06-05"   # last_updated
$infoSheet.Range("B15").Value = 91   # popularity
